# Move the "_GoBack" bookmark (and the space run that preceded it) from the
# end of the "Statistics Tab" Heading 2 paragraph to the end of the
# following "TBD: What statistics will be calculated" paragraph, and add a
# new trailing space run there as well.

$d = $word.ActiveDocument

# Locate the "Statistics Tab" heading paragraph (not the Table of Contents
# entry, which contains many headings concatenated into one paragraph).
$headingPara = $null
$tbdPara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($headingPara -eq $null -and $t.StartsWith("Statistics Tab") -and $t.Length -lt 30) {
        $headingPara = $p
        $tbdPara = $d.Paragraphs.Item($i + 1)
        break
    }
}

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$headingXml = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Heading2"/></w:pPr>' + `
    '<w:bookmarkStart w:id="23" w:name="_Toc5963347"/>' + `
    '<w:r><w:t>Statistics Tab</w:t></w:r>' + `
    '<w:bookmarkEnd w:id="23"/>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '</w:p>'

$tbdXml = '<w:p ' + $wNs + '>' + `
    '<w:r><w:t>TBD: What statistics will be calculated</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:bookmarkStart w:id="24" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="24"/>' + `
    '</w:p>'

$full = $d.Range($headingPara.Range.Start, $tbdPara.Range.End)
$full.InsertXML($headingXml + $tbdXml)
